$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Correct the typo in the confirmation year: 1804 -> 1904 for rows 370 to 887 (column B)
$ws.Range("B370:B887").Value = 1904

# Update the active view/selection state to match the saved workbook
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 864
$ws.Range("A888").Select()
